# Update the cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Column => new value (only cells that actually changed)
$updates = @{
    2  = @{ D = "25.918.26";    E = "  -0.67%  " }
    3  = @{ D = "1.741.73";     E = "  -0.45%  " }
    4  = @{ D = "0.9999";       E = "  +0.02%  " }
    5  = @{ D = "248.87";       E = "  +6.08%  " }
    6  = @{ D = "1.0000";       E = "  +0.00%  " }
    7  = @{ D = "0.5094";       E = "  -3.49%  " }
    8  = @{ D = "0.2722";       E = "  -2.41%  " }
    9  = @{ D = "0.06183" }
    10 = @{ D = "1.750.83";     E = "  +0.12%  " }
    11 = @{ D = "0.07238";      E = "  +0.73%  " }
    12 = @{ D = "15.13";        E = "  -1.52%  " }
    13 = @{ D = "0.6480";       E = "  +0.68%  " }
    14 = @{ D = "4.627" }
    15 = @{ D = "77.64";        E = "  -0.97%  " }
    16 = @{ E = "  +0.03%  " }
    17 = @{ D = "0.9997";       E = "  -0.01%  " }
    18 = @{ D = "25.938.64";    E = "  -0.22%  " }
    19 = @{ D = "11.82";        E = "  +1.16%  " }
    20 = @{ D = "0.000006818";  E = "  +1.50%  " }
    21 = @{ D = "1.964.67";     E = "  -0.21%  " }
    22 = @{ E = "  -0.80%  " }
    23 = @{ D = "8.644";        E = "  -1.09%  " }
    24 = @{ D = "5.383";        E = "  +2.85%  " }
    25 = @{ D = "136.30";       E = "  -1.01%  " }
    26 = @{ E = "  -0.70%  " }
    27 = @{ D = "15.23";        E = "  -0.40%  " }
    28 = @{ D = "1.776";        E = "  -1.44%  " }
    29 = @{ D = "105.43";       E = "  +0.42%  " }
    30 = @{ D = "3.907";        E = "  +2.72%  " }
    31 = @{ D = "0.08228";      E = "  -0.68%  " }
    32 = @{ D = "3.640";        E = "  -0.55%  " }
    33 = @{ D = "0.04681";      E = "  +2.33%  " }
    34 = @{ D = "2.655";        E = "  +0.45%  " }
    35 = @{ D = "0.9967";       E = "  -0.70%  " }
    36 = @{ D = "0.6250";       E = "  -1.42%  " }
    37 = @{ E = "  +0.73%  " }
    38 = @{ E = "  +0.42%  " }
    39 = @{ D = "1.923";        E = "  -1.59%  " }
    40 = @{ D = "0.9997";       E = "  -0.01%  " }
    41 = @{ D = "99.48";        E = "  -0.90%  " }
    42 = @{ D = "0.7611";       E = "  +2.03%  " }
    43 = @{ D = "0.3850";       E = "  -1.79%  " }
    44 = @{ D = "4.999";        E = "  -0.46%  " }
    45 = @{ D = "0.1133";       E = "  -1.05%  " }
    46 = @{ D = "6.293";        E = "  -0.72%  " }
    47 = @{ D = "55.46";        E = "  +2.50%  " }
    48 = @{ D = "0.05234" }
    49 = @{ D = "30.70";        E = "  -1.77%  " }
    50 = @{ B = "Decentraland"; C = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D = "0.3412"; E = "  -1.16%  " }
    51 = @{ B = "EnergySwap";   C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";        D = "7.474";   E = "  -1.96%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
